# Driver class grading-comment rewrite (rows 29, 30, 37 — the
# "CustomerMappingDriver Class" questions and the final "Generic /
# Compilation errors" row of Sheet1).
#
# The author replaced three generic grading remarks with more specific
# ones that call out the ClassCastException behaviour seen while
# running the driver class.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 ("For successfully scanning data from input file") comment.
$ws.Range("F29").Value = "(-8) for not completing else condition, not declaring and initlaizing product and customer objects, not adding them to the inventory"

# Row 30 ("For correct and properly aligned output") comment.
$ws.Range("F30").Value = "(-4) For no output displayed due to ClassCastException while running driver class"

# Row 37 ("Generic" / "Compilation errors if any") comment.
$ws.Range("F37").Value = "(-2.5) For getting ClassCastException while running driver class"

# Leave the selection where the author ended up working.
$null = $ws.Range("F37").Select()
